$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Fill in the previously-missing symbol cells (column B) on rows 1327-1328
$ws.Range("B1327").Value = ":ECBASSETSW"
$ws.Range("B1328").Value = ":ECBASSETSW"

# Step 2: Append new weekly rows 1329-1344, using row 1328 as a fully-styled template
# (this carries over the date style on column A and the string type on column B)
$ws.Range("A1328:G1328").Copy($ws.Range("A1329:G1329"))
$ws.Range("A1329").Value = 45261
$ws.Range("C1329").Value = 7002047000000
$ws.Range("D1329").Value = 7002047000000
$ws.Range("E1329").Value = 7002047000000
$ws.Range("F1329").Value = 7002047000000
$ws.Range("G1329").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1330:G1330"))
$ws.Range("A1330").Value = 45268
$ws.Range("C1330").Value = 6993472000000
$ws.Range("D1330").Value = 6993472000000
$ws.Range("E1330").Value = 6993472000000
$ws.Range("F1330").Value = 6993472000000
$ws.Range("G1330").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1331:G1331"))
$ws.Range("A1331").Value = 45271
$ws.Range("C1331").Value = 6987406000000
$ws.Range("D1331").Value = 6987406000000
$ws.Range("E1331").Value = 6987406000000
$ws.Range("F1331").Value = 6987406000000
$ws.Range("G1331").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1332:G1332"))
$ws.Range("A1332").Value = 45271
$ws.Range("C1332").Value = 6987406000000
$ws.Range("D1332").Value = 6987406000000
$ws.Range("E1332").Value = 6987406000000
$ws.Range("F1332").Value = 6987406000000
$ws.Range("G1332").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1333:G1333"))
$ws.Range("A1333").Value = 45271
$ws.Range("C1333").Value = 6987406000000
$ws.Range("D1333").Value = 6987406000000
$ws.Range("E1333").Value = 6987406000000
$ws.Range("F1333").Value = 6987406000000
$ws.Range("G1333").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1334:G1334"))
$ws.Range("A1334").Value = 45271
$ws.Range("C1334").Value = 6987406000000
$ws.Range("D1334").Value = 6987406000000
$ws.Range("E1334").Value = 6987406000000
$ws.Range("F1334").Value = 6987406000000
$ws.Range("G1334").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1335:G1335"))
$ws.Range("A1335").Value = 45271
$ws.Range("C1335").Value = 6987406000000
$ws.Range("D1335").Value = 6987406000000
$ws.Range("E1335").Value = 6987406000000
$ws.Range("F1335").Value = 6987406000000
$ws.Range("G1335").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1336:G1336"))
$ws.Range("A1336").Value = 45271
$ws.Range("C1336").Value = 6987406000000
$ws.Range("D1336").Value = 6987406000000
$ws.Range("E1336").Value = 6987406000000
$ws.Range("F1336").Value = 6987406000000
$ws.Range("G1336").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1337:G1337"))
$ws.Range("A1337").Value = 45271
$ws.Range("C1337").Value = 6987406000000
$ws.Range("D1337").Value = 6987406000000
$ws.Range("E1337").Value = 6987406000000
$ws.Range("F1337").Value = 6987406000000
$ws.Range("G1337").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1338:G1338"))
$ws.Range("A1338").Value = 45271
$ws.Range("C1338").Value = 6987406000000
$ws.Range("D1338").Value = 6987406000000
$ws.Range("E1338").Value = 6987406000000
$ws.Range("F1338").Value = 6987406000000
$ws.Range("G1338").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1339:G1339"))
$ws.Range("A1339").Value = 45271
$ws.Range("C1339").Value = 6987406000000
$ws.Range("D1339").Value = 6987406000000
$ws.Range("E1339").Value = 6987406000000
$ws.Range("F1339").Value = 6987406000000
$ws.Range("G1339").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1340:G1340"))
$ws.Range("A1340").Value = 45271
$ws.Range("C1340").Value = 6987406000000
$ws.Range("D1340").Value = 6987406000000
$ws.Range("E1340").Value = 6987406000000
$ws.Range("F1340").Value = 6987406000000
$ws.Range("G1340").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1341:G1341"))
$ws.Range("A1341").Value = 45271
$ws.Range("C1341").Value = 6987406000000
$ws.Range("D1341").Value = 6987406000000
$ws.Range("E1341").Value = 6987406000000
$ws.Range("F1341").Value = 6987406000000
$ws.Range("G1341").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1342:G1342"))
$ws.Range("A1342").Value = 45271
$ws.Range("C1342").Value = 6987406000000
$ws.Range("D1342").Value = 6987406000000
$ws.Range("E1342").Value = 6987406000000
$ws.Range("F1342").Value = 6987406000000
$ws.Range("G1342").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1343:G1343"))
$ws.Range("A1343").Value = 45271
$ws.Range("C1343").Value = 6987406000000
$ws.Range("D1343").Value = 6987406000000
$ws.Range("E1343").Value = 6987406000000
$ws.Range("F1343").Value = 6987406000000
$ws.Range("G1343").Value = 0

$ws.Range("A1328:G1328").Copy($ws.Range("A1344:G1344"))
$ws.Range("A1344").Value = 45278
$ws.Range("C1344").Value = 6899165000000
$ws.Range("D1344").Value = 6899165000000
$ws.Range("E1344").Value = 6899165000000
$ws.Range("F1344").Value = 6899165000000
$ws.Range("G1344").Value = 0

# Step 3: Append final row 1345 - same as others but WITHOUT column B (matches source data gap)
$ws.Range("A1328").Copy($ws.Range("A1345"))
$ws.Range("A1345").Value = 45278
$ws.Range("C1345").Value = 6899165000000
$ws.Range("D1345").Value = 6899165000000
$ws.Range("E1345").Value = 6899165000000
$ws.Range("F1345").Value = 6899165000000
$ws.Range("G1345").Value = 0
